# Update automatico via Actualizar 02-05-2021 08-17-35
#
# 1) Refresh the timestamp of the last existing batch of rows (436-449)
#    from 44232.32444751982 -> 44232.32444752315
# 2) Append a brand-new batch of 14 rows (450-463) with the same
#    Name/URL/Status pattern used throughout the sheet, stamped with the
#    new run's timestamp (44232.34550189609), including the per-row
#    hyperlink on column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldStamp = 44232.32444752315
for ($r = 436; $r -le 449; $r++) {
    $ws.Cells.Item($r, 4).Value = $oldStamp
}

$newStamp = 44232.34550189609

# row -> Name (col A), cell text shown in col B, hyperlink Address (+ optional
# SubAddress/fragment, same split MapStore already uses elsewhere in the sheet)
$rows = @(
    @{ Row = 450; Name = "Odoo";              Text = "https://www.dataintelligence-group.com/";              Address = "https://www.dataintelligence-group.com/" },
    @{ Row = 451; Name = "Blackbox";          Text = "https://serviciodashboard.azurewebsites.net/";          Address = "https://serviciodashboard.azurewebsites.net/" },
    @{ Row = 452; Name = "PowerBI";           Text = "https://powerbi.microsoft.com/es-es/";                  Address = "https://powerbi.microsoft.com/es-es/" },
    @{ Row = 453; Name = "Dropbox";           Text = "https://www.dropbox.com/";                              Address = "https://www.dropbox.com/" },
    @{ Row = 454; Name = "Odoo";              Text = "https://dataintelligence.store/";                       Address = "https://dataintelligence.store/" },
    @{ Row = 455; Name = "GEE";               Text = "https://app-data-i.users.earthengine.app/";             Address = "https://app-data-i.users.earthengine.app/" },
    @{ Row = 456; Name = "UtilidadesOdoo";    Text = "https://odooutil.azurewebsites.net/";                   Address = "https://odooutil.azurewebsites.net/" },
    @{ Row = 457; Name = "Filtros Dashboard"; Text = "https://filtradordashboard.azurewebsites.net/";         Address = "https://filtradordashboard.azurewebsites.net/" },
    @{ Row = 458; Name = "MapStore";          Text = "https://ide.dataintelligence-group.com/mapstore/#/";    Address = "https://ide.dataintelligence-group.com/mapstore/"; SubAddress = "/" },
    @{ Row = 459; Name = "GeoServer";         Text = "https://ide.dataintelligence-group.com/geoserver/web/?0"; Address = "https://ide.dataintelligence-group.com/geoserver/web/?0" },
    @{ Row = 460; Name = "Tomcat";            Text = "https://ide.dataintelligence-group.com/";               Address = "https://ide.dataintelligence-group.com/" },
    @{ Row = 461; Name = "Shiny";             Text = "https://rpubs.com/dataintelligence/";                   Address = "https://rpubs.com/dataintelligence/" },
    @{ Row = 462; Name = "Github";            Text = "https://github.com/Sud-Austral/";                       Address = "https://github.com/Sud-Austral/" },
    @{ Row = 463; Name = "EZ Exporter";       Text = "https://ezexporter.highviewapps.com/exports/export-profile/"; Address = "https://ezexporter.highviewapps.com/exports/export-profile/" }
)

foreach ($item in $rows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $item.Name
    $ws.Cells.Item($r, 2).Value = $item.Text
    $ws.Cells.Item($r, 3).Value = "Disponible"

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $newStamp
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $bCell = $ws.Cells.Item($r, 2)
    if ($item.SubAddress) {
        $ws.Hyperlinks.Add($bCell, $item.Address, $item.SubAddress)
    } else {
        $ws.Hyperlinks.Add($bCell, $item.Address)
    }
    # Hyperlinks.Add mints its own style variant; re-apply the sheet's
    # shared "Hyperlink" cell style so the new cells match the rest of
    # the column instead of growing the style table.
    $bCell.Style = "Hyperlink"
}
